$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D9: "Third Manual Column Reduce" List-of-Cols -- insert 'txndb' before 'txndbl' and append 'ebitda'
$ws.Range("D9").Value = "['lco', 'lcox', 'lcoxdr', 'lo', 'loxdr', 'mib', 'mibn', 'ppeveb', 'pstkc', 'pstkl', 'pstkn', 'pstkr', 'tstkc', 'tstkn', 'che', 'cicurr', 'cidergl', 'cimii', 'ciother', 'cipen', 'cisecgl', 'citotal', 'dpc', 'dpvieb', 'dv', 'dvp', 'dvpa', 'ib', 'epspx', 'esopnr', 'esopt', 'ibadj', 'ibc', 'ibcom', 'ibmii', 'recch', 'recco', 'rectr', 'reuna', 'sale', 'spced', 'spceeps', 'cshtr_c', 'dvpsp_c', 'dvpsx_c', 'prcc_c', 'prch_c', 'prcl_c', 'adjex_c', 'acdo', 'aco', 'acodo', 'acox', 'aldo', 'aocidergl', 'aociother', 'aocipen', 'aodo', 'aox', 'ap', 'ceql', 'intc', 'ivaco', 'ivaeq', 'ivao', 'ivch', 'ivst', 'ivstch', 'pncad', 'pncaeps', 'prcad', 'prcaeps', 'xido', 'xidoc', 'ajex', 'ajp', 'cshfd', 'cshi', 'csho', 'cstk', 'cstkcv', 'cstke', 'dclo', 'dcom', 'dcvsr', 'dcvsub', 'dcvt', 'dd', 'dd1', 'dd2', 'dltis', 'dlto', 'dm', 'dn', 'ds', 'dudd', 'fatc', 'fatc', 'fatn', 'fiao', 'fopox', 'intano', 'mrc1', 'mrcta', 'niadj', 'nopio', 'oiadp', 'oibdp', 'oprepsx', 'pnrsho', 'prsho', 'ppent', 'pstkrv', 'txbco', 'txbcof', 'txdba', 'txdbca', 'txdbcl', 'txdc', 'txdi', 'txditc', 'txndba', 'txndb', 'txndbl', 'txo', 'txp', 'txpd', 'txr', 'acctstd', 'am', 'capxv', 'dc', 'diladj', 'do', 'donr', 'emp', 'esub', 'exre', 'lifr', 'mibt', 'prstkc', 'seqo', 'spi', 'cshtr_f', 'dvpsp_f', 'dvpsx_f', 'prcc_f', 'prch_f', 'prcl_f', 'adjex_f', 'ebitda']"

# D10: "Summary Features" List-of-Cols -- append 'txdb', 'txfo', 'txt', 'txw'
$ws.Range("D10").Value = "['drc', 'drlt', 'rea', 'reajo', 'recta', 'txdb', 'txfo', 'txt', 'txw']"

# ResultShape column updates reflecting the revised column counts above
$ws.Range("E9").Value = "(1243, 90)"
$ws.Range("E10").Value = "(1243, 83)"
$ws.Range("E11").Value = "(348, 157)"
$ws.Range("E12").Value = "(348, 162)"
$ws.Range("E17").Value = "(452, 2)"
$ws.Range("E35").Value = "(348, 162)"
$ws.Range("E36").Value = "(348, 169)"
$ws.Range("E37").Value = "(348, 175)"
$ws.Range("E38").Value = "(348, 176)"
$ws.Range("E39").Value = "(348, 179)"
